$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.206652760505676
$ws.Range("B1").Value = 2.391429424285889
$ws.Range("C1").Value = 6.989912509918213
$ws.Range("D1").Value = 2.286118745803833
$ws.Range("E1").Value = 1.173157215118408
